$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Cells.Item(2,2).Value = 1.02
$ws.Cells.Item(2,3).Value = 1.049372342731084
$ws.Cells.Item(2,4).Value = 1.05524637138613
$ws.Cells.Item(2,5).Value = 1.056519709342214
$ws.Cells.Item(2,6).Value = 1.067132573282867
$ws.Cells.Item(2,9).Value = 1.041848971240415
$ws.Cells.Item(2,10).Value = 1.054411138602306
$ws.Cells.Item(2,11).Value = 1.057987356641866
$ws.Cells.Item(2,12).Value = 1.05925719649649
$ws.Cells.Item(2,13).Value = 1.069841256274515
$ws.Cells.Item(2,14).Value = 1.02188240050593
$ws.Cells.Item(3,2).Value = 1.02
$ws.Cells.Item(3,3).Value = 1.050641670357191
$ws.Cells.Item(3,4).Value = 1.056256662867596
$ws.Cells.Item(3,5).Value = 1.057640081975788
$ws.Cells.Item(3,6).Value = 1.068320380919643
$ws.Cells.Item(3,9).Value = 1.042151848842784
$ws.Cells.Item(3,10).Value = 1.055328154213817
$ws.Cells.Item(3,11).Value = 1.058810549063882
$ws.Cells.Item(3,12).Value = 1.060190443672389
$ws.Cells.Item(3,13).Value = 1.070843863620556
$ws.Cells.Item(3,14).Value = 1.022193756721284
$ws.Cells.Item(4,2).Value = 1.02
$ws.Cells.Item(4,3).Value = 1.051462580709762
$ws.Cells.Item(4,4).Value = 1.056909679332007
$ws.Cells.Item(4,5).Value = 1.058364941921949
$ws.Cells.Item(4,6).Value = 1.069088813590013
$ws.Cells.Item(4,9).Value = 1.042345815429295
$ws.Cells.Item(4,10).Value = 1.055920595138128
$ws.Cells.Item(4,11).Value = 1.059341890785107
$ws.Cells.Item(4,12).Value = 1.060793633302113
$ws.Cells.Item(4,13).Value = 1.071491882721988
$ws.Cells.Item(4,14).Value = 1.022394741303565
$ws.Cells.Item(5,2).Value = 1.02
$ws.Cells.Item(5,3).Value = 1.05180759117196
$ws.Cells.Item(5,4).Value = 1.057184038778563
$ws.Cells.Item(5,5).Value = 1.058669651556289
$ws.Cells.Item(5,6).Value = 1.069411826386207
$ws.Cells.Item(5,9).Value = 1.042426876921376
$ws.Cells.Item(5,10).Value = 1.056169436277919
$ws.Cells.Item(5,11).Value = 1.059564952300512
$ws.Cells.Item(5,12).Value = 1.061047051945568
$ws.Cells.Item(5,13).Value = 1.071764135443229
$ws.Cells.Item(5,14).Value = 1.022479119849911
$ws.Cells.Item(6,2).Value = 1.02
$ws.Cells.Item(6,3).Value = 1.051865514142619
$ws.Cells.Item(6,4).Value = 1.057230095081056
$ws.Cells.Item(6,5).Value = 1.058720812428448
$ws.Cells.Item(6,6).Value = 1.069466059573738
$ws.Cells.Item(6,9).Value = 1.042440459247179
$ws.Cells.Item(6,10).Value = 1.056211204882636
$ws.Cells.Item(6,11).Value = 1.059602386899669
$ws.Cells.Item(6,12).Value = 1.061089592575503
$ws.Cells.Item(6,13).Value = 1.07180983768108
$ws.Cells.Item(6,14).Value = 1.022493280626074
$ws.Cells.Item(7,2).Value = 1.02
$ws.Cells.Item(7,3).Value = 1.051467191147157
$ws.Cells.Item(7,4).Value = 1.056913345997635
$ws.Cells.Item(7,5).Value = 1.058369013549555
$ws.Cells.Item(7,6).Value = 1.069093129843456
$ws.Cells.Item(7,9).Value = 1.04234690047055
$ws.Cells.Item(7,10).Value = 1.055923921032243
$ws.Cells.Item(7,11).Value = 1.059344872578486
$ws.Cells.Item(7,12).Value = 1.060797020130335
$ws.Cells.Item(7,13).Value = 1.071495521260117
$ws.Cells.Item(7,14).Value = 1.022395869226634
$ws.Cells.Item(8,2).Value = 1.02
$ws.Cells.Item(8,3).Value = 1.049801408426992
$ws.Cells.Item(8,4).Value = 1.055587951806748
$ws.Cells.Item(8,5).Value = 1.056898365322411
$ws.Cells.Item(8,6).Value = 1.067534032176533
$ws.Cells.Item(8,9).Value = 1.041951747733662
$ws.Cells.Item(8,10).Value = 1.054721241684395
$ws.Cells.Item(8,11).Value = 1.058265832161852
$ws.Cells.Item(8,12).Value = 1.05957273415938
$ws.Cells.Item(8,13).Value = 1.070180245147378
$ws.Cells.Item(8,14).Value = 1.021987725210745
$ws.Cells.Item(9,2).Value = 1.02
$ws.Cells.Item(9,3).Value = 1.046862681270108
$ws.Cells.Item(9,4).Value = 1.053246946863102
$ws.Cells.Item(9,5).Value = 1.054306091130179
$ws.Cells.Item(9,6).Value = 1.064785418154973
$ws.Cells.Item(9,9).Value = 1.041239981307299
$ws.Cells.Item(9,10).Value = 1.052594780422623
$ws.Cells.Item(9,11).Value = 1.056354270621054
$ws.Cells.Item(9,12).Value = 1.057410086125491
$ws.Cells.Item(9,13).Value = 1.067856867539925
$ws.Cells.Item(9,14).Value = 1.02126479922716
$ws.Cells.Item(10,2).Value = 1.02
$ws.Cells.Item(10,3).Value = 1.044901076039247
$ws.Cells.Item(10,4).Value = 1.051682503799863
$ws.Cells.Item(10,5).Value = 1.052577269171685
$ws.Cells.Item(10,6).Value = 1.062952040658439
$ws.Cells.Item(10,9).Value = 1.040755047169783
$ws.Cells.Item(10,10).Value = 1.051172210341495
$ws.Cells.Item(10,11).Value = 1.055072990731549
$ws.Cells.Item(10,12).Value = 1.055964666219294
$ws.Cells.Item(10,13).Value = 1.06630402809712
$ws.Cells.Item(10,14).Value = 1.020780316505537
$ws.Cells.Item(11,2).Value = 1.02
$ws.Cells.Item(11,3).Value = 1.044051059952946
$ws.Cells.Item(11,4).Value = 1.051004168655499
$ws.Cells.Item(11,5).Value = 1.05182849575234
$ws.Cells.Item(11,6).Value = 1.062157914124582
$ws.Cells.Item(11,9).Value = 1.040542585595719
$ws.Cells.Item(11,10).Value = 1.050555030152875
$ws.Cells.Item(11,11).Value = 1.054516526176095
$ws.Cells.Item(11,12).Value = 1.055337895497051
$ws.Cells.Item(11,13).Value = 1.065630679763498
$ws.Cells.Item(11,14).Value = 1.020569922972107
$ws.Cells.Item(12,2).Value = 1.02
$ws.Cells.Item(12,3).Value = 1.043735228449959
$ws.Cells.Item(12,4).Value = 1.050752064649439
$ws.Cells.Item(12,5).Value = 1.051550338535824
$ws.Cells.Item(12,6).Value = 1.0618628986135
$ws.Cells.Item(12,9).Value = 1.040463294297733
$ws.Cells.Item(12,10).Value = 1.050325599930825
$ws.Cells.Item(12,11).Value = 1.05430957908888
$ws.Cells.Item(12,12).Value = 1.055104948654403
$ws.Cells.Item(12,13).Value = 1.065380422211278
$ws.Cells.Item(12,14).Value = 1.020491681293444
$ws.Cells.Item(13,2).Value = 1.02
$ws.Cells.Item(13,3).Value = 1.043802979853231
$ws.Cells.Item(13,4).Value = 1.050806148184226
$ws.Cells.Item(13,5).Value = 1.051610005568125
$ws.Cells.Item(13,6).Value = 1.061926182358122
$ws.Cells.Item(13,9).Value = 1.040480319468792
$ws.Cells.Item(13,10).Value = 1.050374821761967
$ws.Cells.Item(13,11).Value = 1.054353981344134
$ws.Cells.Item(13,12).Value = 1.05515492274127
$ws.Cells.Item(13,13).Value = 1.065434109944513
$ws.Cells.Item(13,14).Value = 1.02050846857904
$ws.Cells.Item(14,2).Value = 1.02
$ws.Cells.Item(14,3).Value = 1.044024955234314
$ws.Cells.Item(14,4).Value = 1.050983332528595
$ws.Cells.Item(14,5).Value = 1.051805503802425
$ws.Cells.Item(14,6).Value = 1.062133528910942
$ws.Cells.Item(14,9).Value = 1.040536038978084
$ws.Cells.Item(14,10).Value = 1.050536069104152
$ws.Cells.Item(14,11).Value = 1.054499425003111
$ws.Cells.Item(14,12).Value = 1.055318642829788
$ws.Cells.Item(14,13).Value = 1.065609996386027
$ws.Cells.Item(14,14).Value = 1.020563457373964
$ws.Cells.Item(15,2).Value = 1.02
$ws.Cells.Item(15,3).Value = 1.04416170854939
$ws.Cells.Item(15,4).Value = 1.051092483029275
$ws.Cells.Item(15,5).Value = 1.051925952754988
$ws.Cells.Item(15,6).Value = 1.062261276392592
$ws.Cells.Item(15,9).Value = 1.040570320077768
$ws.Cells.Item(15,10).Value = 1.050635394732822
$ws.Cells.Item(15,11).Value = 1.054589004281218
$ws.Cells.Item(15,12).Value = 1.055419498062419
$ws.Cells.Item(15,13).Value = 1.065718346421703
$ws.Cells.Item(15,14).Value = 1.020597325555969
$ws.Cells.Item(16,2).Value = 1.02
$ws.Cells.Item(16,3).Value = 1.044957475291523
$ws.Cells.Item(16,4).Value = 1.051727503115215
$ws.Cells.Item(16,5).Value = 1.052626958784371
$ws.Cells.Item(16,6).Value = 1.063004738514177
$ws.Cells.Item(16,9).Value = 1.040769095191713
$ws.Cells.Item(16,10).Value = 1.051213145144528
$ws.Cells.Item(16,11).Value = 1.055109886292837
$ws.Cells.Item(16,12).Value = 1.05600624392096
$ws.Cells.Item(16,13).Value = 1.066348695642926
$ws.Cells.Item(16,14).Value = 1.020794266741878
$ws.Cells.Item(17,2).Value = 1.02
$ws.Cells.Item(17,3).Value = 1.045456468251054
$ws.Cells.Item(17,4).Value = 1.052125586625652
$ws.Cells.Item(17,5).Value = 1.053066631133953
$ws.Cells.Item(17,6).Value = 1.063471020973635
$ws.Cells.Item(17,9).Value = 1.040893116580796
$ws.Cells.Item(17,10).Value = 1.051575230549783
$ws.Cells.Item(17,11).Value = 1.055436175686682
$ws.Cells.Item(17,12).Value = 1.056374053296512
$ws.Cells.Item(17,13).Value = 1.066743839015025
$ws.Cells.Item(17,14).Value = 1.020917639183749
$ws.Cells.Item(18,2).Value = 1.02
$ws.Cells.Item(18,3).Value = 1.045747462027938
$ws.Cells.Item(18,4).Value = 1.052357693262821
$ws.Cells.Item(18,5).Value = 1.053323067148253
$ws.Cells.Item(18,6).Value = 1.0637429705167
$ws.Cells.Item(18,9).Value = 1.040965216742748
$ws.Cells.Item(18,10).Value = 1.051786313394834
$ws.Cells.Item(18,11).Value = 1.055626334535765
$ws.Cells.Item(18,12).Value = 1.056588504120674
$ws.Cells.Item(18,13).Value = 1.066974227102572
$ws.Cells.Item(18,14).Value = 1.020989541493221
$ws.Cells.Item(19,2).Value = 1.02
$ws.Cells.Item(19,3).Value = 1.045846673245157
$ws.Cells.Item(19,4).Value = 1.052436820628851
$ws.Cells.Item(19,5).Value = 1.053410502308017
$ws.Cells.Item(19,6).Value = 1.063835694086946
$ws.Cells.Item(19,9).Value = 1.040989760447439
$ws.Cells.Item(19,10).Value = 1.051858267669645
$ws.Cells.Item(19,11).Value = 1.055691146661254
$ws.Cells.Item(19,12).Value = 1.056661611752808
$ws.Cells.Item(19,13).Value = 1.067052767880715
$ws.Cells.Item(19,14).Value = 1.021014048376075
$ws.Cells.Item(20,2).Value = 1.02
$ws.Cells.Item(20,3).Value = 1.045402937299024
$ws.Cells.Item(20,4).Value = 1.052082885224914
$ws.Cells.Item(20,5).Value = 1.05301946027062
$ws.Cells.Item(20,6).Value = 1.063420995894595
$ws.Cells.Item(20,9).Value = 1.040879835029175
$ws.Cells.Item(20,10).Value = 1.051536394169384
$ws.Cells.Item(20,11).Value = 1.055401184509591
$ws.Cells.Item(20,12).Value = 1.056334599758399
$ws.Cells.Item(20,13).Value = 1.06670145343897
$ws.Cells.Item(20,14).Value = 1.020904408569709
$ws.Cells.Item(21,2).Value = 1.02
$ws.Cells.Item(21,3).Value = 1.043959591740403
$ws.Cells.Item(21,4).Value = 1.05093116003792
$ws.Cells.Item(21,5).Value = 1.051747935274179
$ws.Cells.Item(21,6).Value = 1.062072471684046
$ws.Cells.Item(21,9).Value = 1.040519641290753
$ws.Cells.Item(21,10).Value = 1.050488590821959
$ws.Cells.Item(21,11).Value = 1.05445660242346
$ws.Cells.Item(21,12).Value = 1.055270435122081
$ws.Cells.Item(21,13).Value = 1.065558206259151
$ws.Cells.Item(21,14).Value = 1.020547267093195
$ws.Cells.Item(22,2).Value = 1.02
$ws.Cells.Item(22,3).Value = 1.043051536291534
$ws.Cells.Item(22,4).Value = 1.050206213063989
$ws.Cells.Item(22,5).Value = 1.050948305737939
$ws.Cells.Item(22,6).Value = 1.061224359449618
$ws.Cells.Item(22,9).Value = 1.040291011232716
$ws.Cells.Item(22,10).Value = 1.049828741138811
$ws.Cells.Item(22,11).Value = 1.053861250951212
$ws.Cells.Item(22,12).Value = 1.054600563090014
$ws.Cells.Item(22,13).Value = 1.064838555678801
$ws.Cells.Item(22,14).Value = 1.020322184601082
$ws.Cells.Item(23,2).Value = 1.02
$ws.Cells.Item(23,3).Value = 1.043532968504469
$ws.Cells.Item(23,4).Value = 1.050590598697566
$ws.Cells.Item(23,5).Value = 1.051372221416234
$ws.Cells.Item(23,6).Value = 1.061673983424203
$ws.Cells.Item(23,9).Value = 1.04041241752768
$ws.Cells.Item(23,10).Value = 1.050178640544069
$ws.Cells.Item(23,11).Value = 1.054176996537352
$ws.Cells.Item(23,12).Value = 1.05495575040695
$ws.Cells.Item(23,13).Value = 1.065220136894617
$ws.Cells.Item(23,14).Value = 1.02044155585998
$ws.Cells.Item(24,2).Value = 1.02
$ws.Cells.Item(24,3).Value = 1.045427125818506
$ws.Cells.Item(24,4).Value = 1.05210218042409
$ws.Cells.Item(24,5).Value = 1.0530407748072
$ws.Cells.Item(24,6).Value = 1.063443600152183
$ws.Cells.Item(24,9).Value = 1.040885837130435
$ws.Cells.Item(24,10).Value = 1.051553943014702
$ws.Cells.Item(24,11).Value = 1.055416996011461
$ws.Cells.Item(24,12).Value = 1.056352427380255
$ws.Cells.Item(24,13).Value = 1.066720605941714
$ws.Cells.Item(24,14).Value = 1.020910387096227
$ws.Cells.Item(25,2).Value = 1.02
$ws.Cells.Item(25,3).Value = 1.047622833087187
$ws.Cells.Item(25,4).Value = 1.053852811780939
$ws.Cells.Item(25,5).Value = 1.054976361537221
$ws.Cells.Item(25,6).Value = 1.065496163706886
$ws.Cells.Item(25,9).Value = 1.041425824850538
$ws.Cells.Item(25,10).Value = 1.053145383261053
$ws.Cells.Item(25,11).Value = 1.056849666461761
$ws.Cells.Item(25,12).Value = 1.057969819886519
$ws.Cells.Item(25,13).Value = 1.068458200562215
$ws.Cells.Item(25,14).Value = 1.021452137135813
